# Fill in self-assessment grades on "Group and Self Assessment" sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Group and Self Assessment")
$ws2 = $wb.Worksheets.Item("User Stories")

# Row 10 self/peer assessment grades (team member self-evaluation entries)
$ws1.Range("D10").Value = 4
$ws1.Range("E10").Value = 4
$ws1.Range("F10").Value = 4
$ws1.Range("G10").Value = 3
$ws1.Range("H10").Value = 2
$ws1.Range("I10").Value = 2
$ws1.Range("J10").Value = 4

# Fill in the "User Stories" sheet - student responsible + assessment for US17-US19
$ws2.Range("B12").Value = "1201123 José Silva"
$ws2.Range("C12").Value = 4

$ws2.Range("B13").Value = "1201123 José Silva"
$ws2.Range("C13").Value = 4

$ws2.Range("B14").Value = "1201123 José Silva"
$ws2.Range("C14").Value = 4

# Restore selections to match the final editing state (ws2 selected last so it
# remains the active/visible tab, matching the saved workbook view).
$ws1.Range("J10").Select()
$ws2.Range("C14").Select()
